$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1, 1).Range.Text = "30 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"
$t.Cell(1, 2).Range.Text = "19 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"
$t.Cell(1, 3).Range.Text = "32 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "3|    |" + $nl + "2|    |"
$t.Cell(2, 1).Range.Text = "32 x 56" + $nl + "  5    6" + $nl + "  ----" + $nl + "3|    |" + $nl + "2|    |"
$t.Cell(2, 2).Range.Text = "78 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "7|    |" + $nl + "8|    |"
$t.Cell(2, 3).Range.Text = "49 x 56" + $nl + "  5    6" + $nl + "  ----" + $nl + "4|    |" + $nl + "9|    |"
$t.Cell(3, 1).Range.Text = "14 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "1|    |" + $nl + "4|    |"
$t.Cell(3, 2).Range.Text = "37 x 24" + $nl + "  2    4" + $nl + "  ----" + $nl + "3|    |" + $nl + "7|    |"
$t.Cell(3, 3).Range.Text = "15 x 83" + $nl + "  8    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"
$t.Cell(4, 1).Range.Text = "53 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(4, 2).Range.Text = "67 x 17" + $nl + "  1    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"
$t.Cell(4, 3).Range.Text = "58 x 37" + $nl + "  3    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
$t.Cell(5, 1).Range.Text = "18 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "1|    |" + $nl + "8|    |"
$t.Cell(5, 2).Range.Text = "59 x 50" + $nl + "  5    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
$t.Cell(5, 3).Range.Text = "59 x 11" + $nl + "  1    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
